# Insert a new data row at row 229 of Sheet1, shifting the existing
# rows 229-259 down to 230-260, then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 229..259 down by one row.
$ws.Rows("229:229").Insert()

# Populate the newly inserted row 229 with the new record.
$ws.Cells.Item(229, 1).Value  = 11
$ws.Cells.Item(229, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(229, 3).Value  = "Bíobío"
$ws.Cells.Item(229, 4).Value  = 45154
$ws.Cells.Item(229, 5).Value  = 8
$ws.Cells.Item(229, 6).Value  = 100112032
$ws.Cells.Item(229, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(229, 8).Value  = "Sin especificar"
$ws.Cells.Item(229, 9).Value  = "Primera"
$ws.Cells.Item(229, 10).Value = 250
$ws.Cells.Item(229, 11).Value = 16000
$ws.Cells.Item(229, 12).Value = 17000
$ws.Cells.Item(229, 13).Value = 16600
$ws.Cells.Item(229, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(229, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(229, 16).Value = 332
$ws.Cells.Item(229, 17).Value = 50
$ws.Cells.Item(229, 18).Value = "Hortaliza"
